# Updates "想去人数" (want-to-go count) figures across the workbook's
# sheets to match the refreshed data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 58
$ws1.Range("F7").Value = 1849
$ws1.Range("F8").Value = 3171
$ws1.Range("F24").Value = 4649
$ws1.Range("F28").Value = 5375
$ws1.Range("F30").Value = 1167
$ws1.Range("F31").Value = 230
$ws1.Range("F36").Value = 111
$ws1.Range("F37").Value = 760
$ws1.Range("F40").Value = 699

# --- 本地生活 (Local life) sheet ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 215
$ws3.Range("F3").Value = 1068

# --- 全部类型 (All types) sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 215
$ws4.Range("F4").Value = 1068
$ws4.Range("F9").Value = 58
$ws4.Range("F10").Value = 1849
$ws4.Range("F12").Value = 3171
$ws4.Range("F29").Value = 4649
$ws4.Range("F33").Value = 5375
$ws4.Range("F35").Value = 1167
$ws4.Range("F36").Value = 230
$ws4.Range("F42").Value = 111
$ws4.Range("F43").Value = 760
$ws4.Range("F46").Value = 699
